$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45/46: coin rank swap (ApeXProtocol moves above Fetch.AI)
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"

# Price (D) and Volume(1h) (E) updates
$updates = @{
    2 = @("66.339.16", "  +8.82%  ")
    3 = @("3.462.34", "  +12.73%  ")
    4 = @($null, "  -0.36%  ")
    5 = @("186.05", "  +13.62%  ")
    6 = @("543.77", "  +8.32%  ")
    7 = @("3.446.93", "  +12.33%  ")
    8 = @("0.600", "  +5.01%  ")
    9 = @($null, "  -0.06%  ")
    10 = @("0.624", "  +9.14%  ")
    11 = @("0.148", "  +19.62%  ")
    12 = @("54.08", "  +6.74%  ")
    13 = @("0.0000265", "  +10.96%  ")
    14 = @($null, "  +8.34%  ")
    15 = @("4.020.57", "  +11.72%  ")
    16 = @("3.457.80", "  +11.78%  ")
    17 = @($null, "  +9.15%  ")
    18 = @("66.284.27", "  +8.68%  ")
    19 = @($null, "  +9.73%  ")
    20 = @("11.62", "  +12.15%  ")
    21 = @("0.980", "  +7.16%  ")
    22 = @("411.67", "  +17.63%  ")
    23 = @("84.01", "  +8.43%  ")
    24 = @($null, "  +8.88%  ")
    25 = @($null, "  +11.35%  ")
    26 = @("11.06", "  +5.46%  ")
    27 = @("2.86", "  +15.07%  ")
    28 = @($null, "  +1.08%  ")
    29 = @("11.70", "  +10.26%  ")
    30 = @("8.65", "  +12.11%  ")
    31 = @("29.76", "  +10.49%  ")
    32 = @("650.60", "  +4.83%  ")
    33 = @("6.58", "  +8.81%  ")
    34 = @("11.57", "  +7.87%  ")
    35 = @($null, "  +10.46%  ")
    36 = @("58.92", "  +1.58%  ")
    37 = @("38.03", "  +9.97%  ")
    38 = @("1.00", "  +0.06%  ")
    39 = @("0.0₃0796", "  +23.09%  ")
    40 = @("0.385", "  +8.05%  ")
    41 = @($null, "  +16.31%  ")
    42 = @("3.33", "  +20.79%  ")
    43 = @("0.998", "  -0.49%  ")
    44 = @("2.976.98", "  +8.54%  ")
    47 = @("2.86", "  +18.52%  ")
    48 = @("0.0411", "  +11.40%  ")
    49 = @("2.67", "  +4.49%  ")
    50 = @("8.73", "  +21.30%  ")
    51 = @("0.129", "  +9.18%  ")
    45 = @("3.39", "  +19.50%  ")
    46 = @("2.61", "  +10.32%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $dVal = $pair[0]
    $eVal = $pair[1]
    if ($dVal -ne $null) {
        $ws.Cells.Item($row, 4).Value = "'" + $dVal
        $ws.Cells.Item($row, 4).ClearFormats()
    }
    $ws.Cells.Item($row, 5).Value = "'" + $eVal
    $ws.Cells.Item($row, 5).ClearFormats()
}
